# This script reproduces the effect of inserting a brand-new transaction
# record as row 2 of the worksheet (pushing all existing data rows down
# by one), then filling the new row 2 with the new transaction's data.
#
# Because a blind Rows(...).Insert() in this environment ends up
# re-deriving cell formatting in a way that does not match the original
# workbook's style table, we instead manually shift every data row down
# by one (copying cell-by-cell, from the bottom of the sheet upwards so
# that source data is never overwritten before it is read) and then
# populate the freed-up row 2 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold data for every "normal" transaction row (2-138 before
# the edit / 3-139 after the edit).
$mainCols = @("E", "N", "P", "T")
# Columns that only carry blank, date-formatted placeholder cells, and
# only exist for rows 2-126 before the edit (2-127 after the edit).
$extraCols = @("K", "S", "AB")

$lastExtraSrcRow = 126

# Work from the bottom of the sheet upward so every source row is read
# before it gets overwritten by the row below's shift.
for ($destRow = 139; $destRow -ge 3; $destRow--) {
    $srcRow = $destRow - 1

    foreach ($col in $mainCols) {
        $srcCell = $ws.Range("$col$srcRow")
        $dstCell = $ws.Range("$col$destRow")
        $srcCell.Copy($dstCell)
    }

    if ($srcRow -le $lastExtraSrcRow) {
        foreach ($col in $extraCols) {
            $srcCell = $ws.Range("$col$srcRow")
            $dstCell = $ws.Range("$col$destRow")
            $srcCell.Copy($dstCell)
        }
    }
}

# Populate the new row 2 with the newly-added transaction.
$ws.Range("E2").Value = "Deposit"
$ws.Range("N2").Value = "Crypto"
$ws.Range("P2").Value = "ERC"
$ws.Range("T2").Value = 1060.809

$ws.Range("K2").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("S2").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("AB2").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# Reflect the new selection recorded in the saved workbook.
$ws.Range("T2").Select()
